# B1--and-B2-PowerPoint.pptx edit
#
# 1) The table on slide 5 switches from table style
#    {A6188BA0-2A5C-4DFF-AE81-55D625F6C985} to
#    {A8505DCC-BC66-4893-A6BF-FE477DA7D60D}.
# 2) The deck's applied theme colour palette changes from the "Integral /
#    Red Violet" scheme over to the stock "Office Theme / Office" palette
#    (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink), which is what the
#    before/after theme XML swap amounts to for the theme that is actually
#    applied to the slides.

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on slide 5 -------------------------------------
$slide = $p.Slides.Item(5)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{A8505DCC-BC66-4893-A6BF-FE477DA7D60D}")
    }
}

# --- 2. Swap the theme palette to the stock Office Theme colours ----------
function Set-ThemeRgb($scheme, $index, $r, $g, $b) {
    $scheme.Colors($index).RGB = $r + ($g * 256) + ($b * 65536)
}

$themeColors = $slide.ThemeColorScheme

Set-ThemeRgb $themeColors 1  0x00 0x00 0x00   # dk1
Set-ThemeRgb $themeColors 2  0xFF 0xFF 0xFF   # lt1
Set-ThemeRgb $themeColors 3  0x44 0x54 0x6A   # dk2
Set-ThemeRgb $themeColors 4  0xE7 0xE6 0xE6   # lt2
Set-ThemeRgb $themeColors 5  0x5B 0x9B 0xD5   # accent1
Set-ThemeRgb $themeColors 6  0xED 0x7D 0x31   # accent2
Set-ThemeRgb $themeColors 7  0xA5 0xA5 0xA5   # accent3
Set-ThemeRgb $themeColors 8  0xFF 0xC0 0x00   # accent4
Set-ThemeRgb $themeColors 9  0x44 0x72 0xC4   # accent5
Set-ThemeRgb $themeColors 10 0x70 0xAD 0x47   # accent6
Set-ThemeRgb $themeColors 11 0x05 0x63 0xC1   # hlink
Set-ThemeRgb $themeColors 12 0x95 0x4F 0x72   # folHlink
